# Update results for Steel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Iron & steel / Hydrogen demand value updated
$ws.Range("B3").Value = 16128.74827949199

# Non-metallic minerals column: minor recalculated values
$ws.Range("D6").Value = 4301.342070008422
$ws.Range("D8").Value = 848.7396134232448
